$p = $ppt.ActivePresentation

# Slide 5: "Session-Reflector that supports this TLV, " -> "STAMP Session-Reflector that supports this TLV, "
$s5 = $p.Slides.Item(5)
$shape5 = $s5.Shapes.Item(5)
$para = $shape5.TextFrame.TextRange.Paragraphs(4)
$run = $para.Runs(1)
$run.Text = "STAMP Session-Reflector that supports this TLV, "
# The shape auto-fits to its text (spAutoFit); restore the original height so only
# the text content changes, matching the source edit (a pure text replacement).
$shape5.Height = 260.69633

# Slide 7: "Session-Reflector transmits test packet in-band on the same incoming link in the reverse direction"
#       -> "STAMP Session-Reflector transmits test packet in-band on the same incoming link in the reverse direction"
$s7 = $p.Slides.Item(7)
$shape7 = $s7.Shapes.Item(2)
$para = $shape7.TextFrame.TextRange.Paragraphs(4)
$run = $para.Runs(1)
$run.Text = "STAMP Session-Reflector transmits test packet in-band on the same incoming link in the reverse direction"

# Slide 7: "The Session-Reflector does not transmit reply test packet to the Session-Sender and terminates the Session-Sender test packet"
#       -> "The STAMP Session-Reflector does not transmit reply test packet to the STAMP Session-Sender and terminates the Session-Sender test packet"
$para = $shape7.TextFrame.TextRange.Paragraphs(7)
$run = $para.Runs(1)
$run.Text = "The STAMP Session-Reflector does not transmit reply test packet to the STAMP Session-Sender and terminates the Session-Sender test packet"

# Slide 8: "For this, the Session-Sender can specify in the test packet the receiving destination address for the Session-Reflector reply test packet"
#       -> "STAMP Session-Sender can specify in the test packet the receiving destination address for the STAMP Session-Reflector reply test packet"
$s8 = $p.Slides.Item(8)
$shape8 = $s8.Shapes.Item(2)
$para = $shape8.TextFrame.TextRange.Paragraphs(3)
$run = $para.Runs(1)
$run.Text = "STAMP Session-Sender can specify in the test packet the receiving destination address for the STAMP Session-Reflector reply test packet"

# Slide 9: "For SR path, Session-Reflector reply test packet may need to be sent in-band on a specific return SR path"
#       -> "For an SR path, STAMP Session-Reflector reply test packet may need to be sent in-band on a specific return SR path"
$s9 = $p.Slides.Item(9)
$shape9 = $s9.Shapes.Item(2)
$para = $shape9.TextFrame.TextRange.Paragraphs(1)
$run = $para.Runs(1)
$run.Text = "For an SR path, STAMP Session-Reflector reply test packet may need to be sent in-band on a specific return SR path"

# Slide 9: "Avoid signaling and maintaining dynamic state on Session-Reflector for the return path for each STAMP test session (each session-id, source-address) "
#       -> "Avoid signaling and maintaining dynamic state on STAMP Session-Reflector for the return path for each STAMP test session (each session-id, source-address) "
$para = $shape9.TextFrame.TextRange.Paragraphs(3)
$run = $para.Runs(1)
$run.Text = "Avoid signaling and maintaining dynamic state on STAMP Session-Reflector for the return path for each STAMP test session (each session-id, source-address) "

# Slide 9: "Order of 10K SR Policy (that can have multiple candidate-paths and each candidate-path can have multiple segment-lists)"
#       -> "Can be order of 10K SR Policy (that can also have ECMPs)"
$para = $shape9.TextFrame.TextRange.Paragraphs(4)
$run = $para.Runs(1)
$run.Text = "Can be order of 10K SR Policy (that can also have ECMPs)"
